$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 392, shifting existing rows 392-448 down to 393-449.
$ws.Rows.Item(392).Insert()

# Populate the newly inserted row 392 with a new data record
# (matching the surrounding rows' constant columns and the new date/volume).
$ws.Range("A392").Value = 8
$ws.Range("B392").Value = "Terminal La Palmera de La Serena"
$ws.Range("C392").Value = "Coquimbo"
$ws.Range("D392").Value = 45131
$ws.Range("E392").Value = 4
$ws.Range("F392").Value = 100112012
$ws.Range("G392").Value = "Espinaca"
$ws.Range("H392").Value = "Sin especificar"
$ws.Range("I392").Value = "Primera"
$ws.Range("J392").Value = 1000
$ws.Range("K392").Value = 500
$ws.Range("L392").Value = 600
$ws.Range("M392").Value = 550
$ws.Range("N392").Value = "`$/atado 300 a 500 gramos"
$ws.Range("O392").Value = "Provincia del Elquí"
$ws.Range("P392").Value = 1100
$ws.Range("Q392").Value = 0.5
$ws.Range("R392").Value = "Hortaliza"
